$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 83 so the old data (rows
# 83-101) is pushed down intact to rows 85-103, leaving rows 83-84 free for
# the new week's "Primera"/"Segunda" readings.
$ws.Rows("83:84").Insert()

# Row 83: new "Primera" reading for date 45173 (2023-09-04)
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 45173
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100108
$ws.Range("H83").Value = "Tropicales y subtropicales"
$ws.Range("I83").Value = 100108004
$ws.Range("J83").Value = "Papaya"
$ws.Range("K83").Value = "Cultivar IV Región"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 68
$ws.Range("N83").Value = 18000
$ws.Range("O83").Value = 18000
$ws.Range("P83").Value = 18000
$ws.Range("Q83").Value = "$/bandeja 10 kilos"
$ws.Range("R83").Value = "Provincia del Elquí"
$ws.Range("S83").Value = 1800
$ws.Range("T83").Value = 10

# Row 84: new "Segunda" reading, same date
$ws.Range("A84").Value = 3
$ws.Range("B84").Value = "Femacal de La Calera"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 45173
$ws.Range("E84").Value = 5
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100108
$ws.Range("H84").Value = "Tropicales y subtropicales"
$ws.Range("I84").Value = 100108004
$ws.Range("J84").Value = "Papaya"
$ws.Range("K84").Value = "Cultivar IV Región"
$ws.Range("L84").Value = "Segunda"
$ws.Range("M84").Value = 60
$ws.Range("N84").Value = 15000
$ws.Range("O84").Value = 15000
$ws.Range("P84").Value = 15000
$ws.Range("Q84").Value = "$/bandeja 10 kilos"
$ws.Range("R84").Value = "Provincia del Elquí"
$ws.Range("S84").Value = 1500
$ws.Range("T84").Value = 10
